$wb = $excel.ActiveWorkbook

# --- Sheet 1 "Prix Spot": add new day column BL (16-aug) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy the style from the previous header cell (BK1) onto BL1 so it keeps
# the bold/border header formatting, then set its text.
$ws1.Range("BK1").Copy($ws1.Range("BL1"))
$ws1.Range("BL1").Value = "16-aug"

$ws1.Range("BL2").Value = 76.70999999999999
$ws1.Range("BL3").Value = 74.77
$ws1.Range("BL4").Value = 69.38
$ws1.Range("BL5").Value = 51.61
$ws1.Range("BL6").Value = 44.57
$ws1.Range("BL7").Value = 45.13
$ws1.Range("BL8").Value = 38.97
$ws1.Range("BL9").Value = 48.95
$ws1.Range("BL10").Value = 54.8
$ws1.Range("BL11").Value = 60.93
$ws1.Range("BL12").Value = 26.28
$ws1.Range("BL13").Value = 3.8
$ws1.Range("BL14").Value = 0.17
$ws1.Range("BL15").Value = -0.01
$ws1.Range("BL16").Value = -0.01
$ws1.Range("BL17").Value = -0.01
$ws1.Range("BL18").Value = 3.6
$ws1.Range("BL19").Value = 21.52
$ws1.Range("BL20").Value = 68.8
$ws1.Range("BL21").Value = 86.40000000000001
$ws1.Range("BL22").Value = 83.18000000000001
$ws1.Range("BL23").Value = 93.2
$ws1.Range("BL24").Value = 95.09999999999999
$ws1.Range("BL25").Value = 86.43000000000001

# --- Sheet 2 "Gaz": add new row 61 (2025-08-14) ---
$ws2 = $wb.Worksheets.Item("Gaz")

# Force the date cell to stay plain text (matching the existing rows)
# instead of being auto-converted to an Excel date serial number, then
# drop the resulting formatting so no style index is attached, just like
# the other data rows.
$ws2.Range("A61").NumberFormat = "@"
$ws2.Range("A61").Value = "2025-08-14"
$ws2.Range("A61").ClearFormats()
$ws2.Range("B61").Value = 31.325

# --- Sheet 3 "CO2": add new row 61 (2025-08-14) ---
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Range("A61").NumberFormat = "@"
$ws3.Range("A61").Value = "2025-08-14"
$ws3.Range("A61").ClearFormats()
$ws3.Range("B61").Value = 70.48
